$d = $word.ActiveDocument

# Move to the very end of the document body and add a fresh paragraph
# after the last existing one (the ANTHROPIC_API_KEY= line).
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter() | Out-Null

# Grab the newly created (currently empty) trailing paragraph and
# replace its contents with three distinct runs, matching the target
# OOXML exactly: "GIT_TOKEN=", a single preserved space, then the token.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:r><w:t>GIT_TOKEN=</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:r><w:t>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</w:t></w:r>" +
       "</w:p>"

$paraRange.InsertXML($xml) | Out-Null

Write-Output "GIT_TOKEN paragraph inserted"
